$d = $word.ActiveDocument

$d.Content.Find.Execute("441÷4=110, 1", $false, $false, $false, $false, $false, $true, 1, $false, "444÷4=111, 0", 2) | Out-Null
$d.Content.Find.Execute("958÷9=106, 4", $false, $false, $false, $false, $false, $true, 1, $false, "394÷6=65, 4", 2) | Out-Null
$d.Content.Find.Execute("239÷6=39, 5", $false, $false, $false, $false, $false, $true, 1, $false, "238÷3=79, 1", 2) | Out-Null
$d.Content.Find.Execute("475÷5=95, 0", $false, $false, $false, $false, $false, $true, 1, $false, "941÷5=188, 1", 2) | Out-Null
$d.Content.Find.Execute("302÷8=37, 6", $false, $false, $false, $false, $false, $true, 1, $false, "432÷7=61, 5", 2) | Out-Null
$d.Content.Find.Execute("622÷9=69, 1", $false, $false, $false, $false, $false, $true, 1, $false, "202÷6=33, 4", 2) | Out-Null
$d.Content.Find.Execute("319÷8=39, 7", $false, $false, $false, $false, $false, $true, 1, $false, "680÷8=85, 0", 2) | Out-Null
$d.Content.Find.Execute("605÷2=302, 1", $false, $false, $false, $false, $false, $true, 1, $false, "247÷7=35, 2", 2) | Out-Null
$d.Content.Find.Execute("373÷7=53, 2", $false, $false, $false, $false, $false, $true, 1, $false, "997÷7=142, 3", 2) | Out-Null
$d.Content.Find.Execute("637÷7=91, 0", $false, $false, $false, $false, $false, $true, 1, $false, "801÷6=133, 3", 2) | Out-Null
$d.Content.Find.Execute("709÷8=88, 5", $false, $false, $false, $false, $false, $true, 1, $false, "375÷8=46, 7", 2) | Out-Null
$d.Content.Find.Execute("930÷6=155, 0", $false, $false, $false, $false, $false, $true, 1, $false, "839÷3=279, 2", 2) | Out-Null
$d.Content.Find.Execute("783÷6=130, 3", $false, $false, $false, $false, $false, $true, 1, $false, "863÷9=95, 8", 2) | Out-Null
$d.Content.Find.Execute("685÷8=85, 5", $false, $false, $false, $false, $false, $true, 1, $false, "711÷3=237, 0", 2) | Out-Null
$d.Content.Find.Execute("918÷3=306, 0", $false, $false, $false, $false, $false, $true, 1, $false, "681÷6=113, 3", 2) | Out-Null
$d.Content.Find.Execute("689÷7=98, 3", $false, $false, $false, $false, $false, $true, 1, $false, "401÷8=50, 1", 2) | Out-Null
$d.Content.Find.Execute("371÷3=123, 2", $false, $false, $false, $false, $false, $true, 1, $false, "116÷9=12, 8", 2) | Out-Null
$d.Content.Find.Execute("574÷8=71, 6", $false, $false, $false, $false, $false, $true, 1, $false, "995÷2=497, 1", 2) | Out-Null
$d.Content.Find.Execute("159÷3=53, 0", $false, $false, $false, $false, $false, $true, 1, $false, "305÷9=33, 8", 2) | Out-Null
$d.Content.Find.Execute("931÷6=155, 1", $false, $false, $false, $false, $false, $true, 1, $false, "565÷2=282, 1", 2) | Out-Null
$d.Content.Find.Execute("812÷7=116, 0", $false, $false, $false, $false, $false, $true, 1, $false, "310÷5=62, 0", 2) | Out-Null
$d.Content.Find.Execute("144÷7=20, 4", $false, $false, $false, $false, $false, $true, 1, $false, "390÷7=55, 5", 2) | Out-Null
$d.Content.Find.Execute("238÷2=119, 0", $false, $false, $false, $false, $false, $true, 1, $false, "110÷2=55, 0", 2) | Out-Null
$d.Content.Find.Execute("523÷5=104, 3", $false, $false, $false, $false, $false, $true, 1, $false, "773÷6=128, 5", 2) | Out-Null
$d.Content.Find.Execute("170÷8=21, 2", $false, $false, $false, $false, $false, $true, 1, $false, "536÷3=178, 2", 2) | Out-Null
